$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-05-04 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-05 Monday", 2) | Out-Null

# Update each table cell value directly by row/column to avoid ambiguity
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "80-55=25"
$t.Cell(1, 2).Range.Text = "54+27=81"
$t.Cell(1, 3).Range.Text = "46+21=67"
$t.Cell(1, 4).Range.Text = "19+17=36"
$t.Cell(1, 5).Range.Text = "9+22=31"

$t.Cell(2, 1).Range.Text = "17+73=90"
$t.Cell(2, 2).Range.Text = "53-2=51"
$t.Cell(2, 3).Range.Text = "94-18=76"
$t.Cell(2, 4).Range.Text = "57+20=77"
$t.Cell(2, 5).Range.Text = "83-64=19"

$t.Cell(3, 1).Range.Text = "46+6=52"
$t.Cell(3, 2).Range.Text = "60-36=24"
$t.Cell(3, 3).Range.Text = "54+36=90"
$t.Cell(3, 4).Range.Text = "15-4=11"
$t.Cell(3, 5).Range.Text = "83+11=94"

$t.Cell(4, 1).Range.Text = "83-83=0"
$t.Cell(4, 2).Range.Text = "49+14=63"
$t.Cell(4, 3).Range.Text = "60-2=58"
$t.Cell(4, 4).Range.Text = "15+27=42"
$t.Cell(4, 5).Range.Text = "1+76=77"

$t.Cell(5, 1).Range.Text = "41-30=11"
$t.Cell(5, 2).Range.Text = "95-18=77"
$t.Cell(5, 3).Range.Text = "42+54=96"
$t.Cell(5, 4).Range.Text = "89-47=42"
$t.Cell(5, 5).Range.Text = "29-23=6"

$t.Cell(6, 1).Range.Text = "35-34=1"
$t.Cell(6, 2).Range.Text = "6-1=5"
$t.Cell(6, 3).Range.Text = "54-12=42"
$t.Cell(6, 4).Range.Text = "70-35=35"
$t.Cell(6, 5).Range.Text = "61-53=8"

$t.Cell(7, 1).Range.Text = "33-31=2"
$t.Cell(7, 2).Range.Text = "69+18=87"
$t.Cell(7, 3).Range.Text = "86-33=53"
$t.Cell(7, 4).Range.Text = "94-9=85"
$t.Cell(7, 5).Range.Text = "91-52=39"

$t.Cell(8, 1).Range.Text = "33+7=40"
$t.Cell(8, 2).Range.Text = "62-18=44"
$t.Cell(8, 3).Range.Text = "40+33=73"
$t.Cell(8, 4).Range.Text = "80-42=38"
$t.Cell(8, 5).Range.Text = "4+33=37"

$t.Cell(9, 1).Range.Text = "23+5=28"
$t.Cell(9, 2).Range.Text = "16+33=49"
$t.Cell(9, 3).Range.Text = "69-68=1"
$t.Cell(9, 4).Range.Text = "71+12=83"
$t.Cell(9, 5).Range.Text = "69+22=91"

$t.Cell(10, 1).Range.Text = "81-52=29"
$t.Cell(10, 2).Range.Text = "38+44=82"
$t.Cell(10, 3).Range.Text = "12+66=78"
$t.Cell(10, 4).Range.Text = "62-48=14"
$t.Cell(10, 5).Range.Text = "94-31=63"

$t.Cell(11, 1).Range.Text = "60-4=56"
$t.Cell(11, 2).Range.Text = "50-39=11"
$t.Cell(11, 3).Range.Text = "80+9=89"
$t.Cell(11, 4).Range.Text = "53+4=57"
$t.Cell(11, 5).Range.Text = "74-55=19"

$t.Cell(12, 1).Range.Text = "72-36=36"
$t.Cell(12, 2).Range.Text = "76-70=6"
$t.Cell(12, 3).Range.Text = "31+29=60"
$t.Cell(12, 4).Range.Text = "17-11=6"
$t.Cell(12, 5).Range.Text = "18+6=24"

$t.Cell(13, 1).Range.Text = "98-55=43"
$t.Cell(13, 2).Range.Text = "76+22=98"
$t.Cell(13, 3).Range.Text = "2+17=19"
$t.Cell(13, 4).Range.Text = "92-63=29"
$t.Cell(13, 5).Range.Text = "78-61=17"

$t.Cell(14, 1).Range.Text = "66-24=42"
$t.Cell(14, 2).Range.Text = "28+45=73"
$t.Cell(14, 3).Range.Text = "27+31=58"
$t.Cell(14, 4).Range.Text = "70-47=23"
$t.Cell(14, 5).Range.Text = "56-46=10"

$t.Cell(15, 1).Range.Text = "20+54=74"
$t.Cell(15, 2).Range.Text = "64+18=82"
$t.Cell(15, 3).Range.Text = "56+5=61"
$t.Cell(15, 4).Range.Text = "14+57=71"
$t.Cell(15, 5).Range.Text = "19+40=59"

$t.Cell(16, 1).Range.Text = "16-10=6"
$t.Cell(16, 2).Range.Text = "68-57=11"
$t.Cell(16, 3).Range.Text = "79-77=2"
$t.Cell(16, 4).Range.Text = "70-32=38"
$t.Cell(16, 5).Range.Text = "3+6=9"

$t.Cell(17, 1).Range.Text = "46-31=15"
$t.Cell(17, 2).Range.Text = "39-16=23"
$t.Cell(17, 3).Range.Text = "50-14=36"
$t.Cell(17, 4).Range.Text = "0+3=3"
$t.Cell(17, 5).Range.Text = "68-58=10"

$t.Cell(18, 1).Range.Text = "97-63=34"
$t.Cell(18, 2).Range.Text = "28-10=18"
$t.Cell(18, 3).Range.Text = "74-43=31"
$t.Cell(18, 4).Range.Text = "94-20=74"
$t.Cell(18, 5).Range.Text = "33-8=25"

$t.Cell(19, 1).Range.Text = "48-27=21"
$t.Cell(19, 2).Range.Text = "50+16=66"
$t.Cell(19, 3).Range.Text = "98-80=18"
$t.Cell(19, 4).Range.Text = "60+35=95"
$t.Cell(19, 5).Range.Text = "62-26=36"

$t.Cell(20, 1).Range.Text = "1+20=21"
$t.Cell(20, 2).Range.Text = "86-43=43"
$t.Cell(20, 3).Range.Text = "42-12=30"
$t.Cell(20, 4).Range.Text = "10+68=78"
$t.Cell(20, 5).Range.Text = "72-61=11"
